# Apply updated metrics to Sheet1: every model row now shares the same
# (constant) metric values, and the model_6_8_N labels in column A have
# been reshuffled into a new order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New order of model indices for rows 2..26 (row r -> model_6_8_<n>)
$modelOrder = @(0, 22, 21, 20, 19, 18, 17, 16, 15, 14, 13, 23, 12, 10, 9, 8, 7, 6, 5, 4, 3, 2, 1, 11, 24)

# New constant metric values shared by every row (columns B..Q)
$b = 0.9996522078040857
$c = 0.7171113838424945
$d = 0.9999034524691026
$e = 0.9999999999996098
$f = 0.9999226298236344
$g = 0.0002064644713353352
$h = 0.1679349027030201
$i = [double]"7.508624229984995e-05"
$j = [double]"7.469808810244658e-14"
$k = [double]"3.754312118727402e-05"
$l = 0.004152944209694312
$m = 0.01436887160967538
$n = 1.000203585675657
$o = 0.01498058422176399
$p = 146.9707644240213
$q = 226.1976930404543

$rowVals = @($b, $c, $d, $e, $f, $g, $h, $i, $j, $k, $l, $m, $n, $o, $p, $q)

$arr = New-Object 'object[,]' 25,17
for ($idx = 0; $idx -lt 25; $idx++) {
    $arr[$idx, 0] = "model_6_8_$($modelOrder[$idx])"
    for ($col = 0; $col -lt 16; $col++) {
        $arr[$idx, $col + 1] = $rowVals[$col]
    }
}

$ws.Range("A2:Q26").Value = $arr
